# Daily attendance processing - swap the order of names in the
# "Recorded By" column (G) from "System, <email>" to "<email>, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Text

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val.Split(",")
        if ($parts.Count -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            if ($first -eq "System") {
                $cell.Value = "$second, $first"
            }
        }
    }
}
